# Requisicion form update: the project now pulls from the MySQL relational
# database instead of Firebase, so the "cabecera" (frente/fecha/grupo de
# suministro/lugar de compra) fields were re-filled with the data coming
# back from the new "requisiciones" tables, the proveedor/partida counters
# were corrected, a new line item (Tijeras punta roma) was captured on the
# detail grid, and the stray text value in the NO.- column of the existing
# "Placa de asero reforzado" row was normalized to a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues - used below to drop a text-literal formula down to a plain
# value without Excel re-parsing the literal (needed for cells whose number
# format would otherwise reinterpret a yyyy-mm-dd/numeric-looking string as
# a date or a number instead of leaving it as text).
$xlPasteValues = -4163

# Header / totals counter
$ws.Range("O5").Value = 3

# GRUPO DE SUMINISTRO checkboxes: MATERIALES DE CONSTRUCCION, REFACCIONES and
# PAPELERIA are now marked ("X")
$ws.Range("H8").Value = "X"
$ws.Range("H9").Value = "X"
$ws.Range("H13").Value = "X"

# FRENTE: MAQUINARIA -> ADMINISTRACION
$ws.Range("C9").Value = "ADMINISTRACION"

# LUGAR DE COMPRA checkbox: COMPRA LOCAL is now marked ("X")
$ws.Range("O9").Value = "X"

# FECHA: 2023-06-20 -> 2023-06-22 (cell is date-formatted, so round-trip the
# literal through a text formula + paste-values to keep it stored as text,
# matching the original shared-string/style layout instead of becoming a
# date serial number)
$ws.Range("C11").Formula = "=""2023-06-22"""
$ws.Range("C11").Copy()
$ws.Range("C11").PasteSpecial($xlPasteValues)

# PROVEEDOR on the first detail row: "provedor1" -> "1" (kept as text, same
# paste-values trick so the purely-numeric literal isn't auto-typed as a
# number)
$ws.Range("K20").Formula = "=""1"""
$ws.Range("K20").Copy()
$ws.Range("K20").PasteSpecial($xlPasteValues)

# New line item captured on row 21
$ws.Range("C21").Formula = "=""1234"""
$ws.Range("C21").Copy()
$ws.Range("C21").PasteSpecial($xlPasteValues)

$ws.Range("D21").Value = "Tijeras punta roma"
$ws.Range("I21").Value = "PZA"
$ws.Range("J21").Value = 20

# Existing "Placa de asero reforzado" row: NO.- becomes a real number
$ws.Range("C30").Value = 1256
